# Apply the changes described by the diff to slide 4 of the presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Hunk 1: move/retag the existing "Content Placeholder 3" table (graphicFrame) ---
$existing = $s.Shapes.Item(1)
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Content Placeholder 3") {
        $existing = $sh
    }
}
$existing.Left = 509016 / 12700
$existing.Top = -1089088 / 12700

# --- Hunk 2: add a brand-new table (graphicFrame) holding the feature-comparison data ---
$left   = 838200 / 12700
$top    = 1825625 / 12700
$width  = 8412480 / 12700
$height = 4414520 / 12700

$newTbl = $s.Shapes.AddTable(9, 4, $left, $top, $width, $height)
$newTbl.Name = "Content Placeholder 3"

$table = $newTbl.Table
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $table.Rows.Item($r).Height = 370840 / 12700
}
$newTbl.Left = $left
$newTbl.Top = $top
$newTbl.Width = $width
$newTbl.Height = $height

$data = @(
    @("Fitur", "Jmlh kategori ", "Jumlah image yang ditampilkan", "Rata-rata precision"),
    @("Color histogram", "100", "12", "30.98 %"),
    @("BDIP", "100", "12", "10.56 %"),
    @("BVLC", "100", "12", "10.65 %"),
    @("BDIP + BVLC", "100", "12", "14.5 %"),
    @("Color histogram + BDIP", "100", "12", "31 %"),
    @("Color histogram + BVLC", "100", "12", "31 %"),
    @("Color histogram + BDIP + BVLC", "100", "12", "31%"),
    @("", "", "", "")
)

for ($r = 1; $r -le $data.Count; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Count; $c++) {
        $text = $row[$c - 1]
        if ($text -ne "") {
            $table.Cell($r, $c).Shape.TextFrame.TextRange.Text = $text
        }
    }
}

Write-Output "edit complete"
